$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9573984742164612
$ws.Range("B1").Value = 1.290390014648438
$ws.Range("C1").Value = 2.113709926605225
$ws.Range("D1").Value = 4.451748847961426
$ws.Range("E1").Value = 2.111733913421631
